# Update "PRAZO ACESSO - DIAS" (column E) values for rows 5-12
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = 14
$ws.Range("E6").Value = 14
$ws.Range("E7").Value = 13
$ws.Range("E8").Value = 13
$ws.Range("E9").Value = 12
$ws.Range("E10").Value = 11
$ws.Range("E11").Value = 11
$ws.Range("E12").Value = 10
